# Update the stale/placeholder marketing-authorisation identifier value to the
# real authorisation number "EU/1/97/049/001" wherever it appears (draft POR
# focused edit).
$wb = $excel.ActiveWorkbook

$newId = "EU/1/97/049/001"

$targets = @(
    "AdministrableProductDefinition",
    "ManufacturedItemDefinition",
    "RegulatedAuthorization",
    "Bundle"
)

foreach ($sheetName in $targets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newId
}

# MedicinalProductDefinition!D2 stores the combined
# "<identifier_value>|<hex digest>" string - only the authorisation number
# part changes, the digest after the pipe stays the same.
$ws6 = $wb.Worksheets.Item("MedicinalProductDefinition")
$ws6.Range("D2").Value = "$newId|0xF79CABF272B6A7EEF104DDDA44E82719"

# Refresh the active-cell / selection bookkeeping on the sheets that moved,
# and switch the active tab back to the first sheet.
$ws1 = $wb.Worksheets.Item("AdministrableProductDefinition")
$ws1.Range("C2").Select()

$ws5 = $wb.Worksheets.Item("ManufacturedItemDefinition")
$ws5.Range("C2").Select()

$ws8 = $wb.Worksheets.Item("PackagedProductDefinition")
$ws8.Range("B2").Select()

$ws9 = $wb.Worksheets.Item("RegulatedAuthorization")
$ws9.Range("C2").Select()

$ws1.Activate()
